$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05038766666666666
$ws.Range("M2").Value = 1.309671333333333
$ws.Range("N2").Value = 3.929014
$ws.Range("O2").Value = 0.05806924226264097
$ws.Range("P2").Value = 0.05806924226264098
$ws.Range("Q2").Value = 0.06599128258688888
$ws.Range("R2").Value = 0.593921543282
$ws.Range("S2").Value = 0.05806924226264097
$ws.Range("T2").Value = 0.05806924226264098

$ws.Range("G3").Value = 0.05038766666666666
$ws.Range("O3").Value = 0.3245116581089107
$ws.Range("P3").Value = 0.3245116581089107
$ws.Range("R3").Value = 3.319045630479
$ws.Range("S3").Value = 0.3245116581089107
$ws.Range("T3").Value = 0.3245116581089107

$ws.Range("G4").Value = 0.05038766666666666
$ws.Range("M4").Value = 4.657910333333334
$ws.Range("N4").Value = 13.973731
$ws.Range("O4").Value = 0.2065261082683789
$ws.Range("P4").Value = 0.2065261082683789
$ws.Range("Q4").Value = 0.2347012332392222
$ws.Range("R4").Value = 2.112311099153
$ws.Range("S4").Value = 0.2065261082683789
$ws.Range("T4").Value = 0.2065261082683789

$ws.Range("G5").Value = 0.05038766666666666
$ws.Range("M5").Value = 9.267122333333333
$ws.Range("N5").Value = 27.801367
$ws.Range("O5").Value = 0.4108929913600695
$ws.Range("P5").Value = 0.4108929913600695
$ws.Range("Q5").Value = 0.4669486710912222
$ws.Range("R5").Value = 4.202538039820999
$ws.Range("S5").Value = 0.4108929913600695
$ws.Range("T5").Value = 0.4108929913600695
